$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rate Details")

# --- Number format update: Total Cost / Spot Cost columns switch from
# currency-with-cents (numFmtId 44) to currency-no-cents (new numFmtId 166)
$newFmt = '_("$"* #,##0_);_("$"* \(#,##0\);_("$"* "-"??_);_(@_)'

$ws.Range("B5").NumberFormat = $newFmt
$ws.Range("I8").NumberFormat = $newFmt
$ws.Range("I9").NumberFormat = $newFmt

# --- Clear the hard-coded 0 values out of the "Spots Allocated" input cells
# (H8/H9) so the cells are blank, keeping their existing style/border.
$ws.Range("H8").ClearContents()
$ws.Range("H9").ClearContents()

# --- Fill in the previously-missing helper cells in row 9 (L9/M9), mirroring
# row 8's existing values/format (Currency for L, Comma-no-decimals for M).
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("L9").NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'
$ws.Range("M9").NumberFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# --- Add whole-number data validation on the "Spots Allocated" input range.
$rng = $ws.Range("H8:H9")
$rng.Validation.Add(1, 1, 1, "1", "999999999999999000000")
$rng.Validation.ErrorTitle = "Error"
$rng.Validation.ErrorMessage = "Entry must be a whole number."
$rng.Validation.ShowInput = $false

# --- Window layout tweak recorded in the workbook view.
$win = $wb.Windows.Item(1)
$win.Left = -120
$win.Top = -120
$win.Width = 27375
$win.Height = 16440
